$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "66.931.00"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.878.97"
$ws.Range("E3").Value = "  +3.92%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.38%  "
Set-TextValue $ws.Range("D5") "428.25"
$ws.Range("E5").Value = "  +2.28%  "
Set-TextValue $ws.Range("D6") "132.07"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "3.874.16"
$ws.Range("E7").Value = "  +4.07%  "
Set-TextValue $ws.Range("D8") "0.614"
$ws.Range("E8").Value = "  -5.48%  "
Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  -0.08%  "
Set-TextValue $ws.Range("D10") "0.732"
$ws.Range("E10").Value = "  -4.89%  "
Set-TextValue $ws.Range("D11") "0.168"
$ws.Range("E11").Value = "  -7.26%  "
Set-TextValue $ws.Range("D12") "0.0000363"
$ws.Range("E12").Value = "  -8.85%  "
Set-TextValue $ws.Range("D13") "40.97"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("D14").Value = "4.491.90"
$ws.Range("E14").Value = "  +3.99%  "
Set-TextValue $ws.Range("D15") "10.14"
$ws.Range("E15").Value = "  -3.45%  "
Set-TextValue $ws.Range("D16") "15.69"
$ws.Range("E16").Value = "  +19.08%  "
$ws.Range("D17").Value = "3.882.56"
$ws.Range("E17").Value = "  +4.14%  "
$ws.Range("E18").Value = "  -1.11%  "
Set-TextValue $ws.Range("D19") "19.67"
$ws.Range("E19").Value = "  -5.02%  "
$ws.Range("D20").Value = "67.210.04"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E21").Value = "  -5.76%  "
Set-TextValue $ws.Range("D22") "408.59"
$ws.Range("E22").Value = "  -8.00%  "
Set-TextValue $ws.Range("D23") "14.50"
$ws.Range("E23").Value = "  -11.70%  "
Set-TextValue $ws.Range("D24") "85.33"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("E25").Value = "  -3.64%  "
Set-TextValue $ws.Range("D26") "37.77"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  +11.38%  "
$ws.Range("E28").Value = "  -3.13%  "
Set-TextValue $ws.Range("D29") "9.60"
$ws.Range("E29").Value = "  -5.95%  "
Set-TextValue $ws.Range("D30") "690.48"
$ws.Range("E30").Value = "  +4.65%  "
$ws.Range("E31").Value = "  -1.22%  "
Set-TextValue $ws.Range("D32") "12.49"
$ws.Range("E32").Value = "  -1.75%  "
Set-TextValue $ws.Range("D33") "2.69"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("E35").Value = "  -7.53%  "
Set-TextValue $ws.Range("D36") "38.83"
$ws.Range("E36").Value = "  -7.71%  "
$ws.Range("D37").Value = "0.0₃0811"
$ws.Range("E37").Value = "  +9.28%  "
$ws.Range("E38").Value = "  -0.09%  "
Set-TextValue $ws.Range("D39") "55.37"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E41").Value = "  -6.87%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("E43").Value = "  -9.40%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "27.53"
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D45") "4.52"
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D46") "148.04"
$ws.Range("E46").Value = "  +1.41%  "
Set-TextValue $ws.Range("D47") "2.08"
$ws.Range("E47").Value = "  -2.45%  "
Set-TextValue $ws.Range("D48") "3.27"
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("E51").Value = "  -4.93%  "
